$wb = $excel.ActiveWorkbook

# --- "zh-cn" sheet: rows 4-7 (1b6e9ea3, a5464e12, d2299f97, ee825937) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4:E7").Value = "ht"
$wsZh.Range("H4:H7").Value = "2016-08-26 14:41:35"

# --- "de-de" sheet: rows 4-7 (same file set) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4:E7").Value = "ht"
# "Latest HO Xliff Generate Date" for this file on the Overview sheet, and the
# matching "Latest Handoff Datetime" on de-de, both share the same text value.
$wsDe.Range("H4:H7").Value = "2016-08-26 14:41:39"

# --- "Overview" sheet: G4:G7 (Latest HO Xliff Generate Date for 1b6e9ea3 row) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-26 14:41:39"
